$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06328177979961902
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 7.035200475096886

$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 10.05705581352019

$ws.Range("B4").Value = 0.006876353814593728
$ws.Range("C4").Value = 0.000002220651329265522
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 9.570906078440121

$ws.Range("B5").Value = 0.02258322285507441
$ws.Range("C5").Value = 0.0001537489499301437
$ws.Range("D5").Value = 157.8057217802531
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 164.3098868293287
